# civil and ece dept change
#
# Two staff members have left the ECE department staff list:
#   - Mr. GOPINATH S      (was row 5)
#   - Mr. RAJASEKAR G     (was row 7)
# Remove their rows from the roster, shifting the remaining staff up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom row up so that the row index of the earlier
# row to remove doesn't shift before we get to it.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(5).Delete()

# The sheet used to run down to row 11 (with some left-over formatted
# but empty cells trailing the data). Restore that trailing formatting
# so the sheet still spans down to row 11 as before.
$ws.Rows.Item(10).RowHeight = 18.75

$ws.Range("D9").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Rows.Item(11).RowHeight = 19.5
$excel.CutCopyMode = 0

$ws.Range("E15").Select()
